$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) and 1h volume-change (E) figures for the refreshed crypto snapshot.
# D-column values that look numeric are written with a Text number format first so
# Excel stores them verbatim (e.g. "1.00", keeping trailing zeros) instead of coercing
# them into a Double and losing the original display (matches the source inlineStr cells).

$ws.Range('D2').Value = '76.209.85'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '2.925.52'
$ws.Range('E3').Value = '  +3.30%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '199.77'
$ws.Range('E5').Value = '  +3.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '595.21'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.552'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('E9').Value = '  +2.29%  '
$ws.Range('D10').Value = '2.927.26'
$ws.Range('E10').Value = '  +3.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.453'
$ws.Range('E11').Value = '  +17.18%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.95'
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('D14').Value = '3.466.10'
$ws.Range('E14').Value = '  +3.50%  '
$ws.Range('D15').Value = '76.147.00'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.07'
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000189'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '2.916.41'
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.26'
$ws.Range('E19').Value = '  +7.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.74'
$ws.Range('E20').Value = '  -3.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '372.04'
$ws.Range('E21').Value = '  -2.49%  '
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.34'
$ws.Range('E23').Value = '  +5.48%  '
$ws.Range('E24').Value = '  +2.37%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '3.057.95'
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.29'
$ws.Range('E27').Value = '  +2.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.68'
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000107'
$ws.Range('E29').Value = '  +2.23%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.14'
$ws.Range('E31').Value = '  +5.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.38'
$ws.Range('E32').Value = '  -4.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '493.54'
$ws.Range('E33').Value = '  -4.18%  '
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.91'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.17'
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.390'
$ws.Range('E38').Value = '  +13.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.109'
$ws.Range('E39').Value = '  +24.48%  '
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '179.44'
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.94'
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.65'
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.13'
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.19'
$ws.Range('E47').Value = '  -2.20%  '
$ws.Range('E48').Value = '  +3.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.89'
$ws.Range('E49').Value = '  +4.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.30'
$ws.Range('E50').Value = '  -2.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.84'
$ws.Range('E51').Value = '  +6.73%  '

# Row 41/42: Kaspa and USDe swap position (ranking reorder) with refreshed values.
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.110'
$ws.Range('E42').Value = '  -7.03%  '
